$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $reversed = $parts[($parts.Length - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
